$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force text interpretation so numeric-looking strings (e.g. "1.001")
    # are not auto-converted into numbers/dates by Excel, matching the
    # original inlineStr cell content. Reset the style back to Normal
    # afterwards so no stray number-format style is left on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.183.14"
Set-TextValue $ws.Range("E2") "  -2.30%  "
Set-TextValue $ws.Range("D3") "1.852.82"
Set-TextValue $ws.Range("E3") "  -1.30%  "
Set-TextValue $ws.Range("D4") "1.001"
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "237.50"
Set-TextValue $ws.Range("E5") "  -1.83%  "
Set-TextValue $ws.Range("D6") "0.6798"
Set-TextValue $ws.Range("D7") "1.002"
Set-TextValue $ws.Range("E7") "  +0.09%  "
Set-TextValue $ws.Range("D8") "0.07717"
Set-TextValue $ws.Range("E8") "  +1.01%  "
Set-TextValue $ws.Range("D9") "0.3031"
Set-TextValue $ws.Range("E9") "  -3.85%  "
Set-TextValue $ws.Range("D10") "23.09"
Set-TextValue $ws.Range("E10") "  -6.13%  "
Set-TextValue $ws.Range("D11") "0.08159"
Set-TextValue $ws.Range("E11") "  -0.24%  "
Set-TextValue $ws.Range("B12") "Polygon"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D12") "0.7220"
Set-TextValue $ws.Range("E12") "  -3.31%  "
Set-TextValue $ws.Range("B13") "WrappedEther"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D13") "1.825.05"
Set-TextValue $ws.Range("E13") "  -5.61%  "
Set-TextValue $ws.Range("D14") "5.184"
Set-TextValue $ws.Range("E14") "  -3.13%  "
Set-TextValue $ws.Range("D15") "89.18"
Set-TextValue $ws.Range("E15") "  -3.80%  "
Set-TextValue $ws.Range("D16") "29.200.95"
Set-TextValue $ws.Range("E16") "  -2.40%  "
Set-TextValue $ws.Range("D17") "0.000007812"
Set-TextValue $ws.Range("E17") "  -1.75%  "
Set-TextValue $ws.Range("E18") "  -5.13%  "
Set-TextValue $ws.Range("D19") "13.12"
Set-TextValue $ws.Range("E19") "  -2.54%  "
Set-TextValue $ws.Range("D20") "233.61"
Set-TextValue $ws.Range("E20") "  -5.65%  "
Set-TextValue $ws.Range("D21") "1.002"
Set-TextValue $ws.Range("E21") "  +0.01%  "
Set-TextValue $ws.Range("D22") "2.102.95"
Set-TextValue $ws.Range("E22") "  -1.78%  "
Set-TextValue $ws.Range("E23") "  +0.01%  "
Set-TextValue $ws.Range("D24") "7.434"
Set-TextValue $ws.Range("E24") "  -4.05%  "
Set-TextValue $ws.Range("D25") "161.88"
Set-TextValue $ws.Range("E25") "  -1.49%  "
Set-TextValue $ws.Range("D26") "8.948"
Set-TextValue $ws.Range("E26") "  -3.54%  "
Set-TextValue $ws.Range("D27") "0.1425"
Set-TextValue $ws.Range("E27") "  -6.45%  "
Set-TextValue $ws.Range("D28") "18.04"
Set-TextValue $ws.Range("E28") "  -3.33%  "
Set-TextValue $ws.Range("D29") "1.955"
Set-TextValue $ws.Range("D30") "1.395"
Set-TextValue $ws.Range("E30") "  -3.00%  "
Set-TextValue $ws.Range("D31") "4.518"
Set-TextValue $ws.Range("E31") "  -0.59%  "
Set-TextValue $ws.Range("D32") "1.484"
Set-TextValue $ws.Range("E32") "  -2.76%  "
Set-TextValue $ws.Range("D33") "4.004"
Set-TextValue $ws.Range("E33") "  -4.77%  "
Set-TextValue $ws.Range("D34") "0.05175"
Set-TextValue $ws.Range("E34") "  -4.98%  "
Set-TextValue $ws.Range("E35") "  -4.31%  "
Set-TextValue $ws.Range("B36") "Frax"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D36") "1.027"
Set-TextValue $ws.Range("E36") "  +2.24%  "
Set-TextValue $ws.Range("B37") "ImmutableX"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D37") "0.7014"
Set-TextValue $ws.Range("E37") "  -5.48%  "
Set-TextValue $ws.Range("D38") "2.665"
Set-TextValue $ws.Range("E38") "  -1.19%  "
Set-TextValue $ws.Range("D39") "0.01844"
Set-TextValue $ws.Range("E39") "  -4.59%  "
Set-TextValue $ws.Range("E40") "  -2.46%  "
Set-TextValue $ws.Range("D41") "0.9106"
Set-TextValue $ws.Range("E41") "  +2.88%  "
Set-TextValue $ws.Range("D42") "1.096.79"
Set-TextValue $ws.Range("E42") "  +5.17%  "
Set-TextValue $ws.Range("D43") "5.988"
Set-TextValue $ws.Range("E43") "  -0.19%  "
Set-TextValue $ws.Range("D44") "0.4268"
Set-TextValue $ws.Range("E44") "  -4.66%  "
Set-TextValue $ws.Range("D45") "69.82"
Set-TextValue $ws.Range("E45") "  -2.82%  "
Set-TextValue $ws.Range("D46") "1.002"
Set-TextValue $ws.Range("E46") "  +0.04%  "
Set-TextValue $ws.Range("D47") "102.67"
Set-TextValue $ws.Range("E47") "  -1.44%  "
Set-TextValue $ws.Range("D48") "1.763"
Set-TextValue $ws.Range("E48") "  -3.13%  "
Set-TextValue $ws.Range("D49") "1.998.10"
Set-TextValue $ws.Range("E49") "  -1.51%  "
Set-TextValue $ws.Range("D50") "9.124"
Set-TextValue $ws.Range("E50") "  -6.02%  "
Set-TextValue $ws.Range("D51") "6.881"
Set-TextValue $ws.Range("E51") "  -8.00%  "
